$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (22 and 23) that no longer exist in the target layout
$ws.Rows("22:23").Delete()

# Rewrite the contents of rows 10-21 to match the rebuilt course-sheet content/order
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C10").Value = "9146830 - Danúbia Caporusso Bargos"

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Give students knowledge of the basics of ecology and impact of engineering activities on the environment. Legal and institutional concepts for sustainable development."
$ws.Range("C11").Value = "Give students knowledge of the basics of ecology and impact of engineering activities on the environment. Legal and institutional concepts for sustainable development."

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2020"
$ws.Range("C13").Value = "01/01/2020"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Concepts and Definitions. Environmental issues. Sustainable Development. Environmental performance. Environmental processes. Environmental standard."
$ws.Range("C14").Value = "Concepts and Definitions. Environmental issues. Sustainable Development. Environmental performance. Environmental processes. Environmental standard."

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C15").Value = "9146830 - Danúbia Caporusso Bargos"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "CONCEPTS AND DEFINITIONS. Environmental Engineering. Environment. Environmental pollution. Environmental Critical Components. ENVIRONMENTAL ISSUES. The Subject of Environmental Transformations. Energy and the Environment. Environmental impacts in the three media. Threatened equilibrium. SUSTAINABLE DEVELOPMENT. Basic Concepts. Legal aspects. ENVIRONMENTAL PERFORMANCE. Environmental Monitoring. Scope of Sustainable Development. Definition of indicators. Definition of Sustainable Indicators. Human Development Indicators - HDI. Environmental Sustainability Indicators. Control of Environmental Processes. ENVIRONMENTAL PROCESSES. ETA Process Control. Water in Nature. Characterization of Water. Water Quality Indicators. ENVIRONMENTAL STANDARD. Ordinance 518. CONAMA 20. Ecological disaster."
$ws.Range("C16").Value = "CONCEPTS AND DEFINITIONS. Environmental Engineering. Environment. Environmental pollution. Environmental Critical Components. ENVIRONMENTAL ISSUES. The Subject of Environmental Transformations. Energy and the Environment. Environmental impacts in the three media. Threatened equilibrium. SUSTAINABLE DEVELOPMENT. Basic Concepts. Legal aspects. ENVIRONMENTAL PERFORMANCE. Environmental Monitoring. Scope of Sustainable Development. Definition of indicators. Definition of Sustainable Indicators. Human Development Indicators - HDI. Environmental Sustainability Indicators. Control of Environmental Processes. ENVIRONMENTAL PROCESSES. ETA Process Control. Water in Nature. Characterization of Water. Water Quality Indicators. ENVIRONMENTAL STANDARD. Ordinance 518. CONAMA 20. Ecological disaster."

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C18").Value = "5464150 - Mariana Consiglio Kasemodel"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas com a utilização de recursos de projeções e audiovisual."
$ws.Range("C19").Value = "Aulas expositivas com a utilização de recursos de projeções e audiovisual."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)"
$ws.Range("C20").Value = "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)"

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada"
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada"

# Row 17 now only carries a column-A label (no B/C content in the target layout)
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# Re-apply the row heights from the diff
$ws.Rows("10").RowHeight = 60
$ws.Rows("11").RowHeight = 60
$ws.Rows("13").RowHeight = 60
$ws.Rows("14").RowHeight = 60
$ws.Rows("15").RowHeight = 120
$ws.Rows("16").RowHeight = 120
$ws.Rows("18").RowHeight = 60
$ws.Rows("19").RowHeight = 60
$ws.Rows("20").RowHeight = 60
$ws.Rows("21").RowHeight = 120
